$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.951.44'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '2.363.81'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''302.63'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").Value = '''95.74'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.501'
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").Value = '''0.489'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").Value = '''34.09'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  +3.78%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '''18.35'
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").Value = '''6.72'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '2.731.27'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '2.360.83'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '''0.791'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '42.917.60'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").Value = '''6.25'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = '''68.04'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").Value = '''235.07'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '''2.17'
$ws.Range("E24").Value = '  -5.20%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '''2.43'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '''24.50'
$ws.Range("E27").Value = '  +0.56%  '
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("D29").Value = '''9.30'
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").Value = '''32.01'
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '''5.01'
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '''130.86'
$ws.Range("E34").Value = '  -11.57%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.0714'
$ws.Range("E35").Value = '  +2.10%  '
$ws.Range("D36").Value = '''1.84'
$ws.Range("E36").Value = '  +3.13%  '
$ws.Range("E37").Value = '  +3.69%  '
$ws.Range("D38").Value = '''4.33'
$ws.Range("E38").Value = '  -2.31%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.81'
$ws.Range("E39").Value = '  +3.00%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '''2.28'
$ws.Range("E40").Value = '  -1.95%  '
$ws.Range("D42").Value = '''21.29'
$ws.Range("E42").Value = '  -3.24%  '
$ws.Range("D43").Value = '1.933.69'
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").Value = '''9.15'
$ws.Range("E47").Value = '  -9.14%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.590.50'
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''1.51'
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("D51").Value = '''71.38'
$ws.Range("E51").Value = '  -1.16%  '
